$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Workbook-level: rename "final allocation" -> "women allocation" and add
#    a new "men allocation" sheet right after it.
# ---------------------------------------------------------------------------
$wsStudents = $wb.Worksheets.Item(1)
$wsWomen = $wb.Worksheets.Item(2)
$wsWomen.Name = "women allocation"

$wsMen = $wb.Worksheets.Add($null, $wsWomen)
$wsMen.Name = "men allocation"

# ---------------------------------------------------------------------------
# 2) "Students" sheet: insert a "Category" column at the front and an
#    "Old FLOOR" column between "OLD HOSTEL" and "FLOOR1", then populate the
#    new columns + three new data rows (one more woman, two men).
# ---------------------------------------------------------------------------
$wsStudents.Columns.Item(1).Insert()
$wsStudents.Columns.Item(5).Insert()
$wsStudents.Columns.Item(1).ColumnWidth = 27.5

# Header row
$wsStudents.Range("A1").Value = "Category"
$wsStudents.Range("A1").WrapText = $false
$wsStudents.Range("A1").HorizontalAlignment = -4108
$wsStudents.Range("A1").VerticalAlignment = -4107
$wsStudents.Range("E1").Value = "Old FLOOR"

# Category column for the existing 6 women rows + "Old FLOOR" values derived
# from the old FLOOR1 numbers that used to live in column D before the insert.
$categories = @("Women", "Women", "Women", "Women", "Women", "Women")
$oldFloors = @(3, 3, 3, 2, 3)
$row = 2
foreach ($cat in $categories) {
    $wsStudents.Cells.Item($row, 1).Value = $cat
    $wsStudents.Cells.Item($row, 1).WrapText = $false
    $wsStudents.Cells.Item($row, 1).HorizontalAlignment = -4108
    $wsStudents.Cells.Item($row, 1).VerticalAlignment = -4107
    $row = $row + 1
}
$wsStudents.Range("E2").Value = 3
$wsStudents.Range("E3").Value = 3
$wsStudents.Range("E4").Value = 3
$wsStudents.Range("E5").Value = 2
$wsStudents.Range("E6").Value = 3

# New row 7 - another woman entry. Copy formatting from row 2 first so the
# blank placeholder styling (s=9/10) becomes the normal data styling (s=8).
$wsStudents.Range("B2:K2").Copy($wsStudents.Range("B7:K7"))
$wsStudents.Range("A7").Value = "Women"
$wsStudents.Range("A7").WrapText = $false
$wsStudents.Range("A7").HorizontalAlignment = -4108
$wsStudents.Range("A7").VerticalAlignment = -4107
$wsStudents.Range("B7").Value = "112201045@smail.iitpkd.ac.in"
$wsStudents.Range("C7").Value = 312
$wsStudents.Range("D7").Value = "Saveri Hostel"
$wsStudents.Range("E7").Value = 2
$wsStudents.Range("F7").Value = 1
$wsStudents.Range("G7").Value = "AF"
$wsStudents.Range("H7").Value = 2
$wsStudents.Range("I7").Value = "AF"
$wsStudents.Range("J7").Value = 3
$wsStudents.Range("K7").Value = "AF"

# New row 8 - first man entry.
$wsStudents.Range("B2:K2").Copy($wsStudents.Range("B8:K8"))
$wsStudents.Range("A8").Value = "Men"
$wsStudents.Range("B8").Value = "112201035@smail.iitpkd.ac.in"
$wsStudents.Range("C8").Value = 335
$wsStudents.Range("D8").Value = "Saveri Hostel"
$wsStudents.Range("E8").Value = 3
$wsStudents.Range("F8").Value = 3
$wsStudents.Range("G8").Value = 335
$wsStudents.Range("H8").Value = 2
$wsStudents.Range("I8").Value = "MBS"
$wsStudents.Range("J8").Value = 2
$wsStudents.Range("K8").Value = "SMS"

# New row 9 - second man entry.
$wsStudents.Range("B2:K2").Copy($wsStudents.Range("B9:K9"))
$wsStudents.Range("A9").Value = "Men"
$wsStudents.Range("B9").Value = "112201035@smail.iitpkd.ac.in"
$wsStudents.Range("C9").Value = 325
$wsStudents.Range("D9").Value = "Malhar Hostel"
$wsStudents.Range("E9").Value = 3
$wsStudents.Range("F9").Value = 3
$wsStudents.Range("G9").Value = "MF1S"
$wsStudents.Range("H9").Value = 2
$wsStudents.Range("I9").Value = "MF1S"
$wsStudents.Range("J9").Value = 1
$wsStudents.Range("K9").Value = "MF1S"

$null = $wsStudents.Range("E10").Select()

# ---------------------------------------------------------------------------
# 3) "women allocation" sheet: insert a "Gender" column at the front, fill
#    it in with "Women" for every existing row, fix the room number that
#    changed for 112201030, and append the new 112201045 row.
# ---------------------------------------------------------------------------
$wsWomen.Columns.Item(1).Insert()
$wsWomen.Range("A1").Value = "Gender"
for ($r = 2; $r -le 6; $r++) {
    $wsWomen.Cells.Item($r, 1).Value = "Women"
}
$wsWomen.Range("C6").Value = 123
$wsWomen.Range("A7").Value = "Women"
$wsWomen.Range("B7").Value = "112201045@smail.iitpkd.ac.in"
$wsWomen.Range("C7").Value = 122

# ---------------------------------------------------------------------------
# 4) "men allocation" sheet: brand new sheet with the Gender/Email/Room
#    headers and two rows of data.
# ---------------------------------------------------------------------------
$wsMen.Range("A1").Value = "Gender"
$wsMen.Range("B1").Value = "Email ID"
$wsMen.Range("C1").Value = "Allocated Room"

$wsMen.Range("A2").Value = "Men"
$wsMen.Range("B2").Value = "112201035@smail.iitpkd.ac.in"
$wsMen.Range("C2").Value = 129

$wsMen.Range("A3").Value = "Men"
$wsMen.Range("B3").Value = "112201035@smail.iitpkd.ac.in"
$wsMen.Range("C3").Value = 250

$wsMen.PageSetup.LeftMargin = 54
$wsMen.PageSetup.RightMargin = 54
$wsMen.PageSetup.TopMargin = 72
$wsMen.PageSetup.BottomMargin = 72
$wsMen.PageSetup.HeaderMargin = 36
$wsMen.PageSetup.FooterMargin = 36

$null = $wsStudents.Activate()
